# Edit: move/rename the "Interrupt Collector details" figure block from the
# "Constraints" paragraph up into the interface summary paragraph, splitting
# the single overview picture into two per-slice pictures.

$d = $word.ActiveDocument
$BR = [char]11

# ---------------------------------------------------------------------
# 1) Insert the new "Interrupt Collector details" figure block right after
#    "Interface containing a basic Interrupt Collector block." in the
#    summary paragraph (before the trailing single-space run).
# ---------------------------------------------------------------------
$newPieces = @(
    " ",
    "Interrupt Collector details:",
    " ",
    ".. figure:: hxs/resources/Eccelerators.Library.IP.InterruptCollectorIfc-Slice0.png",
    "   :scale: 50",
    "   ",
    "   Interrupt Collector details slice0",
    "       ",
    ".. figure:: hxs/resources/Eccelerators.Library.IP.InterruptCollectorIfc-Slice1.png",
    "   :scale: 50",
    "   ",
    "   Interrupt Collector details slice1",
    "       "
)
$newInsertText = [string]::Join($BR, $newPieces)

$oldHead = "Interface containing a basic Interrupt Collector block." + $BR + " "
$newHead = "Interface containing a basic Interrupt Collector block." + $BR + $newInsertText + $BR + " "

$rng1 = $d.Content
$found1 = $rng1.Find.Execute($oldHead, $true, $false, $false, $false, $false, $true, 1, $false, $newHead, 2)
if (-not $found1) {
    Write-Host "ERROR: could not find insertion anchor"
}

# ---------------------------------------------------------------------
# 2) Remove the old "Interrupt Collector details" figure block from the
#    "Constraints" paragraph, leaving only the trailing single-space run
#    right after "   register of an UART. ".
# ---------------------------------------------------------------------
$oldPieces = @(
    "   register of an UART. ",
    "   ",
    "Interrupt Collector details:",
    " ",
    ".. figure:: hxs/resources/InterruptCollectorOverview.png",
    "   :scale: 50",
    "   ",
    "   Interrupt Collector details",
    "                                                                 ",
    " "
)
$oldRemoveText = [string]::Join($BR, $oldPieces)
$newRemoveText = "   register of an UART. " + $BR + " "

$rng2 = $d.Content
$found2 = $rng2.Find.Execute($oldRemoveText, $true, $false, $false, $false, $false, $true, 1, $false, $newRemoveText, 2)
if (-not $found2) {
    Write-Host "ERROR: could not find deletion block"
}

Write-Host "done"
